# Update countries & provincias Spain
# - Insert "Republica de Yibuti" as a new country row right after Grecia,
#   shifting "Consejo Danes para los Refugiados" and "Costa de Marfil" down.
# - Refresh COVID-19 counters for several countries.
# - Bump the "datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Shift the three country names down by one row (82,83,84) to make room
#     for the new "Republica de Yibuti" entry, and give it its own (updated)
#     statistics, while the two shifted countries keep their own stats.
$ws.Range("A82").Value = "Republica de Yibuti"
$ws.Range("A83").Value = "Consejo Danes para los Refugiados"
$ws.Range("A84").Value = "Costa de Marfil"

# --- Numeric data refresh (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1730100
$ws.Range("C4").Value = 4825
$ws.Range("D4").Value = 480321
$ws.Range("E4").Value = 1149014
$ws.Range("G4").Value = 193
$ws.Range("H4").Value = 100765

# India (row 13)
$ws.Range("B13").Value = 154369
$ws.Range("C13").Value = 3576
$ws.Range("D13").Value = 65511
$ws.Range("E13").Value = 84477
$ws.Range("G13").Value = 37
$ws.Range("H13").Value = 4381

# Austria (row 44)
$ws.Range("B44").Value = 16591
$ws.Range("C44").Value = 34
$ws.Range("E44").Value = 718

# Moldavia (row 63)
$ws.Range("B63").Value = 7537
$ws.Range("C63").Value = 232
$ws.Range("E63").Value = 3379
$ws.Range("G63").Value = 7
$ws.Range("H63").Value = 274

# Uzbekistan (row 77)
$ws.Range("D77").Value = 2668
$ws.Range("E77").Value = 673

# Republica de Yibuti - new row (row 82)
$ws.Range("B82").Value = 2697
$ws.Range("C82").Value = 229
$ws.Range("D82").Value = 1185
$ws.Range("E82").Value = 1494
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 18

# Consejo Danes para los Refugiados - shifted (row 83)
$ws.Range("B83").Value = 2546
$ws.Range("C83").Value = 143
$ws.Range("D83").Value = 365
$ws.Range("E83").Value = 2113
$ws.Range("H83").Value = 68

# Costa de Marfil - shifted (row 84)
$ws.Range("B84").Value = 2477
$ws.Range("D84").Value = 1286
$ws.Range("E84").Value = 1161
$ws.Range("H84").Value = 30

# Mozambique (row 155)
$ws.Range("B155").Value = 227
$ws.Range("C155").Value = 14
$ws.Range("E155").Value = 155

# Butan (row 189)
$ws.Range("B189").Value = 28
$ws.Range("C189").Value = 1
$ws.Range("E189").Value = 22

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 17:05"
